$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 5 down to row 6 so the new row matches
# the existing style (date format on A/G, boolean display, etc.)
$ws.Range("A5:I5").Copy()
$ws.Range("A6:I6").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A6").Value = 42649.644780092596
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 9852.6
$ws.Range("D6").Value = 9874.82
$ws.Range("E6").Value = 104.82
$ws.Range("F6").Value = 105.290001
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 0.45
$ws.Range("I6").Value = $false

# Widen column F so the longer values added above fit (the source
# workbook auto-sizes this column to ~10.875 characters wide).
$ws.Columns.Item(6).ColumnWidth = 10
